# Actualización automática 2025-07-29 11:25:09
#
# Updates the July sales figure for RIOS CARRION ANGEL BENIGNO /
# CERAMICAS AL COSTO S.A.S. (PORCELANATO group) from 5.76 to 2654.94,
# and propagates the resulting totals/percentages on the dependent
# summary sheets.

$wb = $excel.ActiveWorkbook

$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsCumplimiento = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")

# 1) "VENTAS POR GRUPO": M6 (PORCELANATO column, RIOS CARRION ANGEL BENIGNO / CERAMICAS AL COSTO S.A.S.)
$wsGrupo.Range("M6").Value = 2654.94

# 2) "VENTA MENSUAL": F6 (julio column, same row) and the julio total F24
$wsMensual.Range("F6").Value = 2654.94
$wsMensual.Range("F24").Value = 33382.14

# 3) "CUMPLIMIENTO MENSUAL": PORCELANATO row (16) VENTA / POR CUMPLIR / CUMPLIMIENTO
$wsCumplimiento.Range("D16").Value = 29845.16
$wsCumplimiento.Range("E16").Value = 8911.380000000001
$wsCumplimiento.Range("F16").Value = 0.7700677098626451

# 4) "CUMPLIMIENTO MENSUAL": TOTAL row (19) VENTA / POR CUMPLIR / CUMPLIMIENTO
$wsCumplimiento.Range("D19").Value = 33382.14
$wsCumplimiento.Range("E19").Value = 24840.86386304603
$wsCumplimiento.Range("F19").Value = 0.5733496691191425

# Widen column E on "CUMPLIMIENTO MENSUAL" (POR CUMPLIR), mirroring the
# autofit side-effect of the updated values (22 -> 23 in the stored OOXML
# <col width=".."/>). The COM ColumnWidth property is offset from the raw
# OOXML width by the default cell padding (5/6 of a character), so we
# subtract it here to land exactly on width="23".
$wsCumplimiento.Columns.Item(5).ColumnWidth = 23 - 5/6
